$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "36.679.83"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.98%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.058.40"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  -0.08%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "243.46"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.665"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.36%  "

$ws.Range("E7").Value = "  +0.01%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "54.67"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -7.37%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "59.60"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("E10").Value = "  -3.61%  "

$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("E12").Value = "  -3.05%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.929"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.91%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "14.74"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -4.38%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.359.84"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.93%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "5.45"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -4.59%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.058.66"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "36.605.63"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.11%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "17.14"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -7.15%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "72.05"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -2.20%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0863"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "238.25"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "

$ws.Range("E23").Value = "  -2.53%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  -3.03%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("E27").Value = "  -3.31%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "164.84"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -2.10%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "20.18"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("E31").Value = "  +7.29%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.10"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -7.86%  "

$ws.Range("E33").Value = "  -5.09%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0597"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.66%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("E36").Value = "  -0.71%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0847"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  -1.74%  "

$ws.Range("E39").Value = "  -4.54%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "4.97"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -5.26%  "

$ws.Range("E41").Value = "  -6.47%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.0215"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -3.17%  "

$ws.Range("E43").Value = "  -3.35%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "94.66"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("E45").Value = "  -7.15%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.413.81"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +8.91%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "7.69"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +13.87%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "16.15"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -5.02%  "

$ws.Range("E49").Value = "  +1.87%  "

$ws.Range("E50").Value = "  -4.19%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.247.75"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.11%  "
